# Edit: insert 3 new price records (rows 110-112) into the daily price log,
# pushing the existing records down by 3 rows (110-216 -> 113-219).
# New data corresponds to date serial 44447 (2021-09-08), qualities
# Especial / Primera / Segunda, with their own Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 110; this shifts rows 110:216
# down to 113:219 and keeps all their existing values/formatting intact.
$ws.Range("A110:T112").EntireRow.Insert()

# Columns that are constant across every data row in this sheet; copy them
# from row 113 (the row that used to be row 110 before the insert) so the
# newly inserted rows match the rest of the table exactly.
$constCols = @("A","B","C","E","F","G","H","I","J","K","Q","R","T")

foreach ($col in $constCols) {
    $srcValue = $ws.Range($col + "113").Value2
    $ws.Range($col + "110").Value2 = $srcValue
    $ws.Range($col + "111").Value2 = $srcValue
    $ws.Range($col + "112").Value2 = $srcValue
}

# Make sure the date column keeps the date-time number format used by the
# rest of column D.
$ws.Range("D110:D112").NumberFormat = $ws.Range("D113").NumberFormat

# New row 110: Fecha 2021-09-08, Calidad Especial
$ws.Range("D110").Value2 = 44447
$ws.Range("L110").Value2 = "Especial"
$ws.Range("M110").Value2 = 160
$ws.Range("N110").Value2 = 30000
$ws.Range("O110").Value2 = 31000
$ws.Range("P110").Value2 = 30500
$ws.Range("S110").Value2 = 4357

# New row 111: Fecha 2021-09-08, Calidad Primera
$ws.Range("D111").Value2 = 44447
$ws.Range("L111").Value2 = "Primera"
$ws.Range("M111").Value2 = 240
$ws.Range("N111").Value2 = 25000
$ws.Range("O111").Value2 = 26000
$ws.Range("P111").Value2 = 25500
$ws.Range("S111").Value2 = 3643

# New row 112: Fecha 2021-09-08, Calidad Segunda
$ws.Range("D112").Value2 = 44447
$ws.Range("L112").Value2 = "Segunda"
$ws.Range("M112").Value2 = 200
$ws.Range("N112").Value2 = 20000
$ws.Range("O112").Value2 = 21000
$ws.Range("P112").Value2 = 20500
$ws.Range("S112").Value2 = 2929
